$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regenerated s_val data (filtered save games), rows 2-49, columns B-E and G (sum).
# F (Win) column is left unchanged.
$data = @(
    @(2, 0.1169995834814548, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 2.426980108624251),
    @(3, 0.6545652718822623, 1.626987699542094, 0.7210945179870265, 13.86384647080068, 16.86649396021207),
    @(4, 0.6545652718822623, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 6.038307959104277),
    @(5, 3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027),
    @(6, 3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027),
    @(7, 1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 3.755628166162433),
    @(8, 0.0006075818656279264, 0.3048912486333797, 3.223369029078222, 0.5333859586016987, 4.062253818178927),
    @(9, 0.1169995834814548, 0.04103571897497393, 0.7210945179870265, 0.5333859586016987, 1.412515779045154),
    @(10, 1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 3.755628166162433),
    @(11, 1.445647641019636, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 3.005019366241741),
    @(12, 1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 3.755628166162433),
    @(13, 3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248),
    @(14, 3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027),
    @(15, 1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 4.327115817150455),
    @(16, 0.6545652718822623, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 2.964545797025059),
    @(17, 3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 8.656069925401464),
    @(18, 1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 4.327115817150455),
    @(19, 1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 4.327115817150455),
    @(20, 0.6545652718822623, 1.626987699542094, 0.1496068669990043, 13.86384647080068, 16.29500630922404),
    @(21, 1.445647641019636, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 6.82939032824165),
    @(22, 3.272327238179451, 1.626987699542094, 18.71679738969934, 13.86384647080068, 37.47995879822157),
    @(23, 3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027),
    @(24, 0.01253208636536152, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 1.571903811587466),
    @(25, 3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027),
    @(26, 0.6545652718822623, 1.626987699542094, 3.223369029078222, 13.86384647080068, 19.36876847130326),
    @(27, 1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 3.755628166162433),
    @(28, 3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248),
    @(29, 3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027),
    @(30, 3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027),
    @(31, 3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 8.656069925401464),
    @(32, 0.6545652718822623, 1.626987699542094, 3.223369029078222, 13.86384647080068, 19.36876847130326),
    @(33, 1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 3.755628166162433),
    @(34, 1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 4.327115817150455),
    @(35, 3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027),
    @(36, 3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248),
    @(37, 3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027),
    @(38, 1.445647641019636, 9.983522426115931, 0.1496068669990043, 13.86384647080068, 25.44262340493525),
    @(39, 3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248),
    @(40, 1.445647641019636, 109.9114832445916, 3.223369029078222, 13.86384647080068, 128.4443463854901),
    @(41, 0.6545652718822623, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 2.964545797025059),
    @(42, 1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 4.327115817150455),
    @(43, 1.445647641019636, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 6.82939032824165),
    @(44, 1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 4.327115817150455),
    @(45, 1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 3.755628166162433),
    @(46, 3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248),
    @(47, 1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 4.327115817150455),
    @(48, 3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027),
    @(49, 3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]  # B: TB
    $ws.Cells.Item($r, 3).Value2 = $row[2]  # C: d2S
    $ws.Cells.Item($r, 4).Value2 = $row[3]  # D: K
    $ws.Cells.Item($r, 5).Value2 = $row[4]  # E: IP
    $ws.Cells.Item($r, 7).Value2 = $row[5]  # G: sum
}
